$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 276.8
$ws.Range("I12").Value = 49
$ws.Range("K12").Value = 49
$ws.Range("M12").Value = 121
$ws.Range("H17").Value = 2417.4285
$ws.Range("J17").Value = 2561.4614
$ws.Range("L17").Value = 7684.3842
$ws.Range("N17").Value = -8020.3842
$ws.Range("H43").Value = 1199.6666
$ws.Range("I43").Value = 1200
$ws.Range("K43").Value = 1200
$ws.Range("M43").Value = -1131
$ws.Range("H80").Value = 1440
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 2066.6667
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 6200.000100000001
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -8196.000100000001
$ws.Range("H83").Value = 1440
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 2066.6667
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 18600.0003
$ws.Range("M83").Value = 492
$ws.Range("N83").Value = -28584.0003
$ws.Range("H86").Value = 5314.2856
$ws.Range("I86").Value = 6925
$ws.Range("J86").Value = 3166.6667
$ws.Range("K86").Value = 6925
$ws.Range("L86").Value = 3166.6667
$ws.Range("M86").Value = -5802
$ws.Range("N86").Value = -5412.6667
$ws.Range("H89").Value = 5314.2856
$ws.Range("I89").Value = 6925
$ws.Range("J89").Value = 3166.6667
$ws.Range("K89").Value = 34625
$ws.Range("L89").Value = 15833.3335
$ws.Range("M89").Value = -29009
$ws.Range("N89").Value = -27065.3335
$ws.Range("H98").Value = 2422.2856
$ws.Range("I98").Value = 1487.5
$ws.Range("K98").Value = 1487.5
$ws.Range("M98").Value = 10.5
$ws.Range("H112").Value = 3173.077
$ws.Range("J112").Value = 4000
$ws.Range("L112").Value = 12000
$ws.Range("N112").Value = -14216
$ws.Range("H122").Value = 2422.2856
$ws.Range("I122").Value = 1487.5
$ws.Range("K122").Value = 4462.5
$ws.Range("M122").Value = -2012.5
$ws.Range("H137").Value = 2980.5833
$ws.Range("I137").Value = 2887.2727
$ws.Range("K137").Value = 8661.8181
$ws.Range("M137").Value = -6111.8181
$ws.Range("H138").Value = 2609.5
$ws.Range("I138").Value = 1438.4
$ws.Range("J138").Value = 2999.8667
$ws.Range("K138").Value = 4315.200000000001
$ws.Range("L138").Value = 8999.6001
$ws.Range("M138").Value = 824.7999999999993
$ws.Range("N138").Value = -19279.6001
$ws.Range("H141").Value = 3379.4285
$ws.Range("I141").Value = 3651.0833
$ws.Range("K141").Value = 10953.2499
$ws.Range("M141").Value = -5773.249899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2710.204
$ws.Range("I32").Value = 2506.383
$ws.Range("K32").Value = 2506.383
$ws.Range("M32").Value = -2219.383
$ws.Range("H110").Value = 826.4545000000001
$ws.Range("I110").Value = 682.3333
$ws.Range("K110").Value = 682.3333
$ws.Range("M110").Value = 1362.6667
$ws.Range("H122").Value = 1810.8334
$ws.Range("I122").Value = 1779.7333
$ws.Range("K122").Value = 5339.199900000001
$ws.Range("M122").Value = -2889.199900000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3171
$ws.Range("I20").Value = 3366.1667
$ws.Range("K20").Value = 3366.1667
$ws.Range("M20").Value = -3119.1667
$ws.Range("H99").Value = 4277.091
$ws.Range("I99").Value = 4477.5557
$ws.Range("J99").Value = 3375
$ws.Range("K99").Value = 4477.5557
$ws.Range("L99").Value = 3375
$ws.Range("M99").Value = -2979.5557
$ws.Range("N99").Value = -6371
$ws.Range("H134").Value = 7446.125
$ws.Range("I134").Value = 7446.125
$ws.Range("K134").Value = 22338.375
$ws.Range("M134").Value = -19803.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("N9").Value = 0
$ws.Range("H31").Value = 1856.091
$ws.Range("J31").Value = 3743.75
$ws.Range("L31").Value = 3743.75
$ws.Range("N31").Value = -4333.75
$ws.Range("H34").Value = 1856.091
$ws.Range("J34").Value = 3743.75
$ws.Range("L34").Value = 3743.75
$ws.Range("N34").Value = -4147.75
$ws.Range("H50").Value = 28216.6
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31250
$ws.Range("H60").Value = 21273.25
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022
$ws.Range("H62").Value = 3749.5
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876
$ws.Range("H65").Value = 3749.5
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380
$ws.Range("H99").Value = 3314.3333
$ws.Range("I99").Value = 4371.5
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 4371.5
$ws.Range("L99").Value = 1200
$ws.Range("M99").Value = -2873.5
$ws.Range("N99").Value = -4196
$ws.Range("H105").Value = 399.5
$ws.Range("I105").Value = 399
$ws.Range("K105").Value = 399
$ws.Range("M105").Value = 1348
$ws.Range("H126").Value = 3314.3333
$ws.Range("I126").Value = 4371.5
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 13114.5
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -10644.5
$ws.Range("N126").Value = -8540
$ws.Range("H140").Value = 92499.5
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2430.2273
$ws.Range("I55").Value = 798.75
$ws.Range("J55").Value = 2792.7778
$ws.Range("K55").Value = 2396.25
$ws.Range("L55").Value = 8378.3334
$ws.Range("M55").Value = -2219.25
$ws.Range("N55").Value = -8732.3334
$ws.Range("H56").Value = 12010
$ws.Range("I56").Value = 12010
$ws.Range("K56").Value = 12010
$ws.Range("M56").Value = -11480
$ws.Range("H103").Value = 438
$ws.Range("I103").Value = 407
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 1221
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = -342
$ws.Range("N103").Value = -3258

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1350
$ws.Range("I113").Value = 1350
$ws.Range("K113").Value = 1350
$ws.Range("M113").Value = 820
$ws.Range("H122").Value = 2635.6875
$ws.Range("I122").Value = 1410
$ws.Range("K122").Value = 4230
$ws.Range("M122").Value = -1780
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = 0

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 48359.3
$ws.Range("I7").Value = 48177
$ws.Range("J7").Value = 50000
$ws.Range("K7").Value = 48177
$ws.Range("L7").Value = 50000
$ws.Range("M7").Value = -48065
$ws.Range("N7").Value = -50224
$ws.Range("H20").Value = 8000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H22").Value = 722.36365
$ws.Range("I22").Value = 986.5
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = 986.5
$ws.Range("L22").Value = 18
$ws.Range("M22").Value = -691.5
$ws.Range("N22").Value = -608
$ws.Range("H27").Value = 722.36365
$ws.Range("I27").Value = 986.5
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = 986.5
$ws.Range("L27").Value = 18
$ws.Range("M27").Value = -879.5
$ws.Range("N27").Value = -232
$ws.Range("H40").Value = 2997.6667
$ws.Range("I40").Value = 2997.6667
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2997.6667
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2861.6667
$ws.Range("H61").Value = 9141.429
$ws.Range("I61").Value = 8495.5
$ws.Range("J61").Value = 10002.667
$ws.Range("K61").Value = 8495.5
$ws.Range("L61").Value = 10002.667
$ws.Range("M61").Value = -8293.5
$ws.Range("N61").Value = -10406.667
$ws.Range("H101").Value = 10347.333
$ws.Range("J101").Value = 10347.333
$ws.Range("L101").Value = 10347.333
$ws.Range("N101").Value = -16837.333
$ws.Range("H113").Value = 9141.429
$ws.Range("I113").Value = 8495.5
$ws.Range("J113").Value = 10002.667
$ws.Range("K113").Value = 8495.5
$ws.Range("L113").Value = 10002.667
$ws.Range("M113").Value = -6325.5
$ws.Range("N113").Value = -14342.667
$ws.Range("H126").Value = 48359.3
$ws.Range("I126").Value = 48177
$ws.Range("J126").Value = 50000
$ws.Range("K126").Value = 144531
$ws.Range("L126").Value = 150000
$ws.Range("M126").Value = -142061
$ws.Range("N126").Value = -154940

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 35799.668
$ws.Range("J41").Value = 35799.668
$ws.Range("L41").Value = 35799.668
$ws.Range("N41").Value = -36579.668
$ws.Range("H81").Value = 3997
$ws.Range("I81").Value = 3997
$ws.Range("K81").Value = 7994
$ws.Range("M81").Value = -6933
$ws.Range("H84").Value = 3997
$ws.Range("I84").Value = 3997
$ws.Range("K84").Value = 39970
$ws.Range("M84").Value = -34666
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0
$ws.Range("H126").Value = 1833
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 1577.2222
$ws.Range("I132").Value = 1423.75
$ws.Range("K132").Value = 4271.25
$ws.Range("M132").Value = -1741.25
